$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E16:E26) so the periods run in ascending
# order (2402 .. 2412) instead of descending (2412 .. 2402).
$periodos = @("2402", "2403", "2404", "2405", "2406", "2407", "2408", "2409", "2410", "2411", "2412")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Update "Valor Mora" column: swap the values that were out of sequence
# so the first row (2402) carries 52000 and the last row (2412) carries
# 64794.
$ws.Range("F16").Value = 52000
$ws.Range("F26").Value = 64794
